$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-01-14 Sunday" "2024-01-15 Monday"
Replace-Text "58÷8=7, 2" "20÷8=2, 4"
Replace-Text "51÷7=7, 2" "36÷4=9, 0"
Replace-Text "73÷2=36, 1" "95÷9=10, 5"
Replace-Text "58÷6=9, 4" "90÷2=45, 0"
Replace-Text "57÷9=6, 3" "77÷4=19, 1"
Replace-Text "99÷5=19, 4" "19÷8=2, 3"
Replace-Text "55÷3=18, 1" "26÷6=4, 2"
Replace-Text "13÷2=6, 1" "79÷5=15, 4"
Replace-Text "26÷9=2, 8" "17÷9=1, 8"
Replace-Text "38÷3=12, 2" "10÷5=2, 0"
Replace-Text "72÷2=36, 0" "77÷5=15, 2"
Replace-Text "85÷4=21, 1" "45÷5=9, 0"
Replace-Text "69÷2=34, 1" "64÷4=16, 0"
Replace-Text "69÷7=9, 6" "62÷9=6, 8"
Replace-Text "85÷9=9, 4" "71÷4=17, 3"
Replace-Text "48÷9=5, 3" "95÷8=11, 7"
Replace-Text "84÷2=42, 0" "19÷8=2, 3"
Replace-Text "34÷3=11, 1" "98÷3=32, 2"
Replace-Text "23÷3=7, 2" "35÷9=3, 8"
Replace-Text "98÷6=16, 2" "39÷6=6, 3"
Replace-Text "35÷3=11, 2" "46÷7=6, 4"
Replace-Text "61÷7=8, 5" "18÷9=2, 0"
Replace-Text "85÷5=17, 0" "84÷7=12, 0"
Replace-Text "10÷4=2, 2" "10÷5=2, 0"
Replace-Text "10÷7=1, 3" "28÷8=3, 4"
